# Weekly fruit/hortaliza data update:
# A new price observation (dated 45166) is inserted as the new first record
# of the "Berenjena" block (row 472), pushing the existing records
# (previously rows 472-512) down by one row (now rows 473-513).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 472, shifting rows 472:512 down to 473:513
# (and the sheet dimension from A1:R512 to A1:R513).
$ws.Rows("472:472").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A472").Value = 3
$ws.Range("B472").Value = "Femacal de La Calera"
$ws.Range("C472").Value = "Coquimbo"
$ws.Range("D472").Value = 45166
$ws.Range("E472").Value = 5
$ws.Range("F472").Value = 100112001
$ws.Range("G472").Value = "Berenjena"
$ws.Range("H472").Value = "Sin especificar"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 110
$ws.Range("K472").Value = 7000
$ws.Range("L472").Value = 7500
$ws.Range("M472").Value = 7318
$ws.Range("N472").Value = "`$/caja 60 unidades"
$ws.Range("O472").Value = "Región de Arica y Parinacota"
$ws.Range("P472").Value = 122
$ws.Range("Q472").Value = 60
$ws.Range("R472").Value = "Hortaliza"
